$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the old "_GoBack" bookmark (it will be re-created later at
#    the position of the last edit, matching Word's own behaviour).
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. Insert the new sentence after "...the INSTALL target."
# ------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("the INSTALL target.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)
$rng.InsertAfter(" You must have write permission to the installation destination or the installation will fail.")
$insertEnd = $rng.End

# ------------------------------------------------------------------
# 3. Re-create the "_GoBack" bookmark right at the end of the text we
#    just inserted (collapsed / zero-length, as Word leaves it after a
#    plain text edit). Adding a bookmark exactly one position before a
#    paragraph mark mis-anchors it, so nudge the paragraph mark out of
#    the way with a throw-away character, add the bookmark, then strip
#    the throw-away character back out.
# ------------------------------------------------------------------
$tempRange = $d.Range($insertEnd, $insertEnd)
$tempRange.InsertAfter("X")
$anchorRange = $d.Range($insertEnd, $insertEnd)
$d.Bookmarks.Add("_GoBack", $anchorRange)
$d.Range($insertEnd, $insertEnd + 1).Delete()

# ------------------------------------------------------------------
# 4. Merge the split runs describing the shared-relwithdebinfo step.
# ------------------------------------------------------------------
$old3 = "Buil" + "d and install the shared-relwithdebinfo" + " configuration"
$new3 = "Build and install the shared-relwithdebinfo configuration"
$d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2)

$old4 = "Repeat steps 2 through 4, setting CMAKE_INSTALL_PREFIX to a different location, for instance VTK_7.0.0-shared" + "-relwithdebinfo" + ". " + "Choose the RelWithDebInfo" + " configuration in Visual Studio and build the install target."
$new4 = "Repeat steps 2 through 4, setting CMAKE_INSTALL_PREFIX to a different location, for instance VTK_7.0.0-shared-relwithdebinfo. Choose the RelWithDebInfo configuration in Visual Studio and build the install target."
$d.Content.Find.Execute($old4, $true, $false, $false, $false, $false, $true, 1, $false, $new4, 2)
